# Insert a new data row at row 172, shifting the existing rows 172:264 down
# to 173:265 (dimension grows from A1:R264 to A1:R265), then populate the
# newly inserted row 172 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(172).Insert()

$ws.Range("A172").Value = 9
$ws.Range("B172").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C172").Value = 'Metropolitana'
$ws.Range("D172").Value = 44806
$ws.Range("D172").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = 100112026
$ws.Range("G172").Value = 'Haba'
$ws.Range("H172").Value = 'Sin especificar'
$ws.Range("I172").Value = 'Primera'
$ws.Range("J172").Value = 68
$ws.Range("K172").Value = 12000
$ws.Range("L172").Value = 13000
$ws.Range("M172").Value = 12412
$ws.Range("N172").Value = '$/saco 25 kilos'
$ws.Range("O172").Value = 'Provincia de Limarí'
$ws.Range("P172").Value = 496
$ws.Range("Q172").Value = 25
$ws.Range("R172").Value = 'Hortaliza'
